$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text by pre-formatting the cell as Text
# before assigning the value (matches original inlineStr / shared-string text cells).

# Row 2
$ws.Range("D2").Value = "98.968.19"
$ws.Range("E2").Value = "  +1.40%  "

# Row 3
$ws.Range("D3").Value = "3.390.97"
$ws.Range("E3").Value = "  +8.60%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "261.91"
$ws.Range("E5").Value = "  +8.36%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "634.43"
$ws.Range("E6").Value = "  +3.64%  "

# Row 7
$ws.Range("E7").Value = "  +24.85%  "

# Row 8
$ws.Range("E8").Value = "  +2.68%  "

# Row 9
$ws.Range("E9").Value = "  -0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.883"
$ws.Range("E10").Value = "  +11.78%  "

# Row 11
$ws.Range("D11").Value = "3.387.01"
$ws.Range("E11").Value = "  +8.51%  "

# Row 12
$ws.Range("E12").Value = "  +1.37%  "

# Row 13
$ws.Range("D13").Value = "98.548.94"
$ws.Range("E13").Value = "  +1.71%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.43"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000250"
$ws.Range("E15").Value = "  +3.24%  "

# Row 16
$ws.Range("D16").Value = "4.005.37"
$ws.Range("E16").Value = "  +8.11%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.56"
$ws.Range("E17").Value = "  +3.46%  "

# Row 18
$ws.Range("D18").Value = "3.386.81"
$ws.Range("E18").Value = "  +8.46%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.60"
$ws.Range("E19").Value = "  -0.06%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.26"
$ws.Range("E20").Value = "  +4.09%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "497.86"
$ws.Range("E21").Value = "  -4.43%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.23"
$ws.Range("E22").Value = "  +9.39%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000214"
$ws.Range("E23").Value = "  +9.95%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.45"
$ws.Range("E24").Value = "  +6.52%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.75"
$ws.Range("E25").Value = "  +3.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.42"
$ws.Range("E26").Value = "  +3.28%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("E27").Value = "  +3.40%  "

# Row 28
$ws.Range("D28").Value = "3.559.64"
$ws.Range("E28").Value = "  +8.21%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.283"
$ws.Range("E29").Value = "  +19.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.201"
$ws.Range("E30").Value = "  +14.92%  "

# Row 32
$ws.Range("E32").Value = "  +6.11%  "

# Row 33
$ws.Range("E33").Value = "  +18.22%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.63"
$ws.Range("E34").Value = "  +6.22%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.05"
$ws.Range("E35").Value = "  +4.76%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.41"
$ws.Range("E36").Value = "  +0.42%  "

# Row 37
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("E37").Value = "  -0.18%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.99"
$ws.Range("E38").Value = "  +6.12%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.474"
$ws.Range("E39").Value = "  +8.09%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "505.41"
$ws.Range("E40").Value = "  +2.10%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "24.87"
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("E43").Value = "  +3.41%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.39"
$ws.Range("E44").Value = "  +5.90%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.793"
$ws.Range("E45").Value = "  +13.95%  "

# Row 46
$ws.Range("E46").Value = "  -0.01%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.85"
$ws.Range("E47").Value = "  -0.41%  "

# Row 48
$ws.Range("E48").Value = "  +1.61%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.848"
$ws.Range("E49").Value = "  +15.70%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.71"
$ws.Range("E50").Value = "  +7.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "46.55"
$ws.Range("E51").Value = "  +4.72%  "
